# Fruta / hortaliza, semanal
# Insert two new weekly observation rows (new rows 210 and 211) into the
# "Hortaliza, Terminal Hortofrutícola Agro Chillán - Repollo" sheet,
# pushing the existing rows 210-254 down to 212-256.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 210 (old row 210 becomes new row 212).
$ws.Rows("210:211").Insert()

# ---- New row 210 ----
$ws.Cells.Item(210, 1).Value2 = 7
$ws.Cells.Item(210, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(210, 3).Value2 = "Ñuble"
$ws.Cells.Item(210, 4).Value2 = 44782
$ws.Cells.Item(210, 5).Value2 = 16
$ws.Cells.Item(210, 6).Value2 = 100112006
$ws.Cells.Item(210, 7).Value2 = "Repollo"
$ws.Cells.Item(210, 8).Value2 = "Crespo record"
$ws.Cells.Item(210, 9).Value2 = "Primera"
$ws.Cells.Item(210, 10).Value2 = 240
$ws.Cells.Item(210, 11).Value2 = 1100
$ws.Cells.Item(210, 12).Value2 = 1300
$ws.Cells.Item(210, 13).Value2 = 1200
$ws.Cells.Item(210, 14).Value2 = "$/unidad"
$ws.Cells.Item(210, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(210, 16).Value2 = 1200
$ws.Cells.Item(210, 17).Value2 = 1
$ws.Cells.Item(210, 18).Value2 = "Hortaliza"

# ---- New row 211 ----
$ws.Cells.Item(211, 1).Value2 = 7
$ws.Cells.Item(211, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(211, 3).Value2 = "Ñuble"
$ws.Cells.Item(211, 4).Value2 = 44782
$ws.Cells.Item(211, 5).Value2 = 16
$ws.Cells.Item(211, 6).Value2 = 100112006
$ws.Cells.Item(211, 7).Value2 = "Repollo"
$ws.Cells.Item(211, 8).Value2 = "Crespo record"
$ws.Cells.Item(211, 9).Value2 = "Segunda"
$ws.Cells.Item(211, 10).Value2 = 200
$ws.Cells.Item(211, 11).Value2 = 900
$ws.Cells.Item(211, 12).Value2 = 900
$ws.Cells.Item(211, 13).Value2 = 900
$ws.Cells.Item(211, 14).Value2 = "$/unidad"
$ws.Cells.Item(211, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(211, 16).Value2 = 900
$ws.Cells.Item(211, 17).Value2 = 1
$ws.Cells.Item(211, 18).Value2 = "Hortaliza"
